# "logger class is added" - the test-automation framework now logs a
# "close" step at the end of the Test Steps sheet and marks every step's
# Status as Executed; the expected login-user-name in Test Data is
# corrected too.

$wb = $excel.ActiveWorkbook
$wsSteps = $wb.Worksheets.Item("Test Steps")
$wsData  = $wb.Worksheets.Item("Test Data")

# --- Test Data: fix the expected login user name ---
$wsData.Range("D3").Value = "Jakay M"

# --- Test Steps: the last existing step now also reports Executed ---
$wsSteps.Range("F7").Value = "Executed"

# --- Test Steps: append a new "close" step, logged as Executed ---
$wsSteps.Range("A7:F7").Copy($wsSteps.Range("A8:F8")) | Out-Null
$wsSteps.Range("A8").Value = "TC_001_Validate_Login_Page"
$wsSteps.Range("B8").Value = ""
$wsSteps.Range("C8").Value = "close"
$wsSteps.Range("D8").Value = ""
$wsSteps.Range("E8").Value = ""
$wsSteps.Range("F8").Value = "Executed"

# --- restore/update view state: Test Data no longer the active tab,
#     Test Steps becomes active with selection at A10 ---
$wsData.Range("D3").Select() | Out-Null
$wsSteps.Range("A10").Select() | Out-Null
